$wb = $excel.ActiveWorkbook

# --- 1. Fix the May sheet: B3/C3/D3 were stored as text ("100","15","100"),
#        convert them to real numeric values (100, 15, 100). ---
$mayWs = $wb.Worksheets.Item("May")
$mayWs.Range("B3").Value = 100
$mayWs.Range("C3").Value = 15
$mayWs.Range("D3").Value = 100

# --- 2. Add a new "June" worksheet after "May" ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$juneWs = $wb.Worksheets.Add($null, $afterSheet)
$juneWs.Name = "June"

# Header row (bold/centered style matching the other sheets' s="1")
$headers = @("Customer Name", "Unit Price", "Consumption Period", "Usage (%)", "Consumption Duration", "Net Price", "Remarks", "Month")
for ($c = 1; $c -le $headers.Length; $c++) {
    $juneWs.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$mayWs.Range("A1:H1").Copy()
$juneWs.Range("A1:H1").PasteSpecial(-4122)

# Row 2 - John
$juneWs.Range("A2").Value = "John"
$juneWs.Range("B2").Value = 100
$juneWs.Range("C2").Value = 15
$juneWs.Range("D2").Value = 50
$juneWs.Range("E2").Value = 0.5
$juneWs.Range("F2").Value = 25
$juneWs.Range("G2").Value = "ok"
$juneWs.Range("H2").Value = "June"

# Row 3 - Jane
$juneWs.Range("A3").Value = "Jane"
$juneWs.Range("B3").Value = 100
$juneWs.Range("C3").Value = 15
$juneWs.Range("D3").Value = 100
$juneWs.Range("E3").Value = 0.5
$juneWs.Range("F3").Value = 50
$juneWs.Range("G3").Value = "ok"
$juneWs.Range("H3").Value = "June"

# Row 4 - Jack (Remarks left blank)
$juneWs.Range("A4").Value = "Jack"
$juneWs.Range("B4").Value = 180
$juneWs.Range("C4").Value = 14
$juneWs.Range("D4").Value = 65
$juneWs.Range("E4").Value = 0.47
$juneWs.Range("F4").Value = 54.98999999999999
$juneWs.Range("H4").Value = "June"

$juneWs.Range("A1").Select()
